# Insert 9 new historical rows (2019-11-18 .. 2019-11-28) for stock 5268 (IKHMAS)
# right before the existing 2019-11-29 row (current row 1054), shifting the rest
# of the table down by 9 rows (A1:I1126 -> A1:I1135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 9 blank rows at 1054, pushing rows 1054:1126 down to 1063:1135.
$ws.Rows("1054:1062").Insert()

# Keep the text columns (date/id/name) stored as text, matching the rest of the
# sheet, instead of letting Excel auto-coerce "5268" to a number or the date
# string to a real date serial.
$ws.Range("B1054:D1062").NumberFormat = "@"

$newRows = @(
    @{ Row=1054; A=1574035200; B="2019-11-18"; E=0.1;                F=0.1;                G=0.09;               H=0.095;              I=6806200 },
    @{ Row=1055; A=1574121600; B="2019-11-19"; E=0.09;               F=0.095;              G=0.08500000000000001;H=0.09;               I=3071300 },
    @{ Row=1056; A=1574208000; B="2019-11-20"; E=0.09;               F=0.095;              G=0.09;               H=0.09;               I=1357500 },
    @{ Row=1057; A=1574294400; B="2019-11-21"; E=0.09;               F=0.09;               G=0.08500000000000001;H=0.09;               I=4187100 },
    @{ Row=1058; A=1574380800; B="2019-11-22"; E=0.09;               F=0.09;               G=0.08500000000000001;H=0.08500000000000001;I=3629700 },
    @{ Row=1059; A=1574640000; B="2019-11-25"; E=0.08500000000000001;F=0.08500000000000001;G=0.08;               H=0.08;               I=9960600 },
    @{ Row=1060; A=1574726400; B="2019-11-26"; E=0.08;               F=0.09;               G=0.08;               H=0.09;               I=7463600 },
    @{ Row=1061; A=1574812800; B="2019-11-27"; E=0.09;               F=0.09;               G=0.08500000000000001;H=0.08500000000000001;I=2291100 },
    @{ Row=1062; A=1574899200; B="2019-11-28"; E=0.09;               F=0.09;               G=0.08;               H=0.08500000000000001;I=2106500 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = "5268"
    $ws.Range("D$row").Value = "IKHMAS"
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
}
